# Updated export json logic

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the worksheet ---
$ws.Name = "Travelrep CZ"

# --- 2. Update the Print_Area defined name to reflect new sheet name / range ---
$printAreaName = $wb.Names.Item(1)
$printAreaName.RefersTo = "='Travelrep CZ'!`$A`$1:`$F`$65"

# --- 3. Move the "Naklady" (Costs) label from row 25 to row 30 ---
# Copy the "Naklady" cell's formatting (style "2") down onto B30 first ...
$ws.Range("B25").Copy()
$ws.Range("B30").PasteSpecial(-4122)
$ws.Range("B30").Value = "Náklady"
# ... then reset B25 back to the plain style (style "1", like its neighbours)
# and clear its text.
$ws.Range("A25").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$ws.Range("B25").Value = ""

# --- 4. Fix small text typos/labels in row 32 and row 44 ---
$ws.Range("E32").Value = "Plné zah."
$ws.Range("F32").Value = "Celk. den"
$ws.Range("C44").Value = "Kapesné:"

# --- 5. Rebuild the bottom section of the sheet (rows 51-65) ---
# Save old values before we start overwriting ranges
# (Value2 is used for reading since it returns the raw value reliably)
$oldE53 = $ws.Range("E53").Value2
$oldA54 = $ws.Range("A54").Value2
$oldA55 = $ws.Range("A55").Value2
$oldB55 = $ws.Range("B55").Value2
$oldE55 = $ws.Range("E55").Value2
$oldF55 = $ws.Range("F55").Value2
$oldE61 = $ws.Range("E61").Value2
$oldA63 = $ws.Range("A63").Value2
$oldE63 = $ws.Range("E63").Value2
$oldA64 = $ws.Range("A64").Value2
$oldC64 = $ws.Range("C64").Value2
$oldE64 = $ws.Range("E64").Value2
$oldE65 = $ws.Range("E65").Value2

# Clear out the old sparse rows 53-65 first
$ws.Range("A53:F65").Clear()

# Apply consistent formatting (copied from row 50, which already has the
# desired style) and row height to every row from 51 to 65
$ws.Range("A50:F50").Copy()
for ($r = 51; $r -le 65; $r++) {
    $ws.Range("A" + $r + ":F" + $r).PasteSpecial(-4122)
    $ws.Rows.Item($r).RowHeight = 12
}

# Row 53
$ws.Range("E53").Value = $oldE53
# Row 54
$ws.Range("A54").Value = $oldA54
# Row 55
$ws.Range("A55").Value = $oldA55
$ws.Range("B55").Value = $oldB55
$ws.Range("E55").Value = $oldE55
$ws.Range("F55").Value = $oldF55
# Row 61
$ws.Range("E61").Value = $oldE61
# Row 63
$ws.Range("A63").Value = $oldA63
$ws.Range("E63").Value = $oldE63
# Row 64
$ws.Range("A64").Value = $oldA64
$ws.Range("C64").Value = $oldC64
$ws.Range("E64").Value = $oldE64
# Row 65
$ws.Range("E65").Value = $oldE65
